# Daily attendance processing - 2025-10-07 09:42:23
# Reorders the "Recorded By" names in column G (System moved to the front)
# and refreshes the dependent attendance count / percentage cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> new "Recorded By" value (column G)
$gUpdates = @{
    2   = "system, System, backup@backdoor.com"
    3   = "System, dnasr281@gmail.com"
    4   = "System, backup@backdoor.com"
    5   = "System, backup@backdoor.com"
    6   = "System, dnasr281@gmail.com"
    11  = "System, dnasr281@gmail.com"
    12  = "System, dnasr281@gmail.com"
    13  = "System, dnasr281@gmail.com"
    14  = "System, dnasr281@gmail.com"
    15  = "System, dnasr281@gmail.com"
    29  = "system, System, backup@backdoor.com"
    30  = "System, dnasr281@gmail.com"
    32  = "System, backup@backdoor.com"
    33  = "System, dnasr281@gmail.com"
    38  = "System, dnasr281@gmail.com"
    39  = "System, dnasr281@gmail.com"
    40  = "System, dnasr281@gmail.com"
    41  = "System, dnasr281@gmail.com"
    42  = "System, dnasr281@gmail.com"
    56  = "system, System, backup@backdoor.com"
    57  = "System, dnasr281@gmail.com"
    58  = "System, backup@backdoor.com"
    59  = "System, backup@backdoor.com"
    60  = "System, dnasr281@gmail.com"
    65  = "System, dnasr281@gmail.com"
    66  = "System, dnasr281@gmail.com"
    67  = "System, dnasr281@gmail.com"
    68  = "System, dnasr281@gmail.com"
    69  = "System, dnasr281@gmail.com"
    84  = "System, backup@backdoor.com"
    85  = "System, backup@backdoor.com"
    86  = "System, dnasr281@gmail.com"
    89  = "System, dnasr281@gmail.com"
    93  = "System, dnasr281@gmail.com"
    110 = "System, backup@backdoor.com"
    111 = "System, backup@backdoor.com"
    112 = "System, dnasr281@gmail.com"
    115 = "System, dnasr281@gmail.com"
    119 = "System, dnasr281@gmail.com"
    136 = "System, backup@backdoor.com"
    137 = "System, backup@backdoor.com"
    138 = "System, dnasr281@gmail.com"
    141 = "System, dnasr281@gmail.com"
    145 = "System, dnasr281@gmail.com"
}

foreach ($row in $gUpdates.Keys) {
    $ws.Range("G$row").Value = $gUpdates[$row]
}

# Attendance count refreshed for row 86 (an extra attendee recorded by System)
$ws.Range("H86").Value = "37/56"

# Cumulative attendance percentage summary cell recalculated.
# Leading apostrophe forces Excel to store this as literal text (matching
# the existing text-formatted percentage cells in this table) instead of
# auto-converting it to a numeric percentage value.
$ws.Range("S18").Value = "'69.1%"
